# added 4wk low sales check
# Updates forecast figures (MyForecast, Inventory Coverage, Stockout Risk,
# Reorder Urgency, Seasonality Index) on the "Forecast Comparison" sheet and
# the rolled-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet --------------------------------------------
# Columns: D = MyForecast, H = Inventory Coverage, I = Stockout Risk,
#          J = Reorder Urgency, L = Seasonality Index

# Row 2 (W10)
$wsForecast.Range("D2").Value = 31
$wsForecast.Range("H2").Value = 10.03
$wsForecast.Range("L2").Value = 0.84

# Row 3 (W11)
$wsForecast.Range("D3").Value = 35
$wsForecast.Range("H3").Value = 8.050000000000001
$wsForecast.Range("L3").Value = 1.04

# Row 4 (W12)
$wsForecast.Range("D4").Value = 35
$wsForecast.Range("H4").Value = 6.89
$wsForecast.Range("L4").Value = 1.08

# Row 5 (W13)
$wsForecast.Range("D5").Value = 34
$wsForecast.Range("H5").Value = 6.21
$wsForecast.Range("L5").Value = 1.04

# Row 6 (W14)
$wsForecast.Range("D6").Value = 32
$wsForecast.Range("H6").Value = 5.43
$wsForecast.Range("L6").Value = 0.97

# Row 7 (W15)
$wsForecast.Range("D7").Value = 33
$wsForecast.Range("H7").Value = 4.3
$wsForecast.Range("L7").Value = 1.05

# Row 8 (W16)
$wsForecast.Range("D8").Value = 36
$wsForecast.Range("H8").Value = 3.03
$wsForecast.Range("L8").Value = 0.84

# Row 9 (W17)
$wsForecast.Range("D9").Value = 39
$wsForecast.Range("H9").Value = 1.9
$wsForecast.Range("I9").Value = "Low"
$wsForecast.Range("J9").Value = "Normal"
$wsForecast.Range("L9").Value = 0.9

# Row 10 (W18)
$wsForecast.Range("D10").Value = 39
$wsForecast.Range("H10").Value = 0.9
$wsForecast.Range("I10").Value = "Low"
$wsForecast.Range("L10").Value = 0.99

# Row 11 (W19)
$wsForecast.Range("D11").Value = 36
$wsForecast.Range("L11").Value = 0.82

# Row 12 (W20)
$wsForecast.Range("D12").Value = 35
$wsForecast.Range("L12").Value = 1.18

# Row 13 (W21)
$wsForecast.Range("D13").Value = 36
$wsForecast.Range("L13").Value = 1.19

# Row 14 (W22)
$wsForecast.Range("D14").Value = 38
$wsForecast.Range("L14").Value = 1.1

# Row 15 (W23)
$wsForecast.Range("D15").Value = 40
$wsForecast.Range("L15").Value = 0.93

# Row 16 (W24)
$wsForecast.Range("D16").Value = 38
$wsForecast.Range("L16").Value = 1.04

# Row 17 (W25)
$wsForecast.Range("D17").Value = 34
$wsForecast.Range("L17").Value = 1.12

# --- Summary sheet ----------------------------------------------------------
# These totals are stored as text (not numbers) in this sheet, just like the
# other Summary values, so a leading apostrophe forces Excel to keep them as
# text instead of auto-converting the numeric-looking string to a number.
$wsSummary.Range("B9").Value = "'576"
$wsSummary.Range("B10").Value = "'278"
$wsSummary.Range("B11").Value = "'136"
$wsSummary.Range("B12").Value = "'40"
$wsSummary.Range("B14").Value = "'31"
